$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix formatting anomaly on existing rows 228:229 (column E style 8 -> 5) ---
# Copy the border/format already used on column A (style index 5) onto E228:E229
$ws.Range("A228").Copy() | Out-Null
$ws.Range("E228:E229").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Add three new rows of data (230-232) ---

# Row 230: Port of Helsinki - South Harbour Live
$ws.Range("A230").Value = "LIVE, SEA, HARBOR"
$ws.Range("B230").Value = "60.16418180859401, 24.959608925070114"
$ws.Range("C230").Value = "Port of Helsinki - South Harbour Live (Live Camera Axis Q6155-E)"
$ws.Range("D230").Value = "Helsinki"
$ws.Range("E230").Value = "Finland"
$ws.Range("F230").Value = "CvOB-Is_yYU"

# Row 231: Port of Helsinki - West harbour - south cam
$ws.Range("A231").Value = "LIVE, SEA, HARBOR"
$ws.Range("B231").Value = "60.1497202443619, 24.914443359739366"
$ws.Range("C231").Value = "Port of Helsinki - West harbour - south cam"
$ws.Range("D231").Value = "Helsinki"
$ws.Range("E231").Value = "Finland"
$ws.Range("F231").Value = "6hPWq2IG08M"

# Row 232: Port of Helsinki - West harbour - north cam
$ws.Range("A232").Value = "LIVE, TRAFFIC"
$ws.Range("B232").Value = "60.15402791777396, 24.918840108266455"
$ws.Range("C232").Value = "Port of Helsinki - West harbour - north cam"
$ws.Range("D232").Value = "Helsinki"
$ws.Range("E232").Value = "Finland"
$ws.Range("F232").Value = "JnJhFYhIjFs"

# Apply the same border formatting (style index 5 / applyFill variant 8) used by
# the neighbouring rows to the new A and E cells of the new rows.
$ws.Range("A229").Copy() | Out-Null
$ws.Range("A230:A232").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E229").Copy() | Out-Null
$ws.Range("E230:E232").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Update the active selection to reflect where the user ended up editing ---
$ws.Range("A233").Select() | Out-Null
